# Update the division-problem worksheet table: each data row (1, 5, 9, 13, 17)
# gets its five cell values replaced with the new set of problems.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("16÷2=", "16÷7=", "24÷9=", "30÷6=", "75÷2=")
    5  = @("58÷2=", "94÷9=", "28÷7=", "28÷6=", "21÷2=")
    9  = @("19÷3=", "10÷4=", "28÷8=", "54÷4=", "76÷7=")
    13 = @("39÷2=", "56÷6=", "85÷9=", "70÷8=", "28÷2=")
    17 = @("94÷8=", "96÷2=", "78÷6=", "69÷4=", "82÷6=")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
